# QA overhaul: fix HTML template, diversify scripts, improve content quality
# Updates a handful of message texts across the ChaylaGreyJourney, NRWaves
# and boosters sheets.

$wb = $excel.ActiveWorkbook

# --- ChaylaGreyJourney ---------------------------------------------------
$ws = $wb.Worksheets.Item("ChaylaGreyJourney")

$ws.Range("B6").Value = "I wanna finish with you papi... I'm so close right now"
$ws.Range("B10").Value = "give me a sec"
$ws.Range("B14").Value = "omg"
$ws.Range("B17").Value = "tell me what you want me to do right now papi... I'll do anything you say 😏"
$ws.Range("B19").Value = "okay you're actually doing things to me right now and it's 100% your fault"
$ws.Range("B20").Value = "did you watch it?"
$ws.Range("B24").Value = "haha I knew you'd like that. you have good taste babe"
$ws.Range("B25").Value = "sooo?"
$ws.Range("B28").Value = "ugh stop I'm literally blushing and I never blush. you're dangerous"
$ws.Range("B32").Value = "so tell me about you, what do you like to do? I wanna know everything"
$ws.Range("B33").Value = "niceee! I'm in Dallas right now, originally from Brazil. I miss the beaches but I love it here too"
$ws.Range("B34").Value = "aww that's so sweet haha. so where are you from?"
$ws.Range("B35").Value = "heyyy papi 💕 thanks for subscribing! you just made my day honestly. what brought you here?"

# --- NRWaves --------------------------------------------------------------
$ws = $wb.Worksheets.Item("NRWaves")

$ws.Range("B2").Value = "I've got something that's going to blow your mind when you get back 🔥"
$ws.Range("B3").Value = "hey babe, don't be a stranger 💕"
$ws.Range("B4").Value = "your loss... this was your exclusive"
$ws.Range("B5").Value = "you're really going to miss out on what I just recorded... 🔥"
$ws.Range("B6").Value = "yo 😏"

# --- boosters ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("boosters")

$ws.Range("B6").Value = "you have no idea what you're doing to me"
